# Append a new paragraph "I am harsha" after the existing "Hello world"
# paragraph (the last paragraph in the document body, just before the
# section properties). The new paragraph inherits the same paragraph
# formatting (ListBullet style, numPr ilvl=0/numId=0, left/hanging indent)
# from the paragraph it follows.

$d = $word.ActiveDocument

$lastPara = $d.Paragraphs.Last
$r = $lastPara.Range
$r.Collapse(0)          # wdCollapseEnd
$r.InsertParagraphAfter()

$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "I am harsha"
